$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 156456
$ws.Range("C4").Value = 147548
$ws.Range("C7").Value = 5.69
$ws.Range("C8").Value = 63.79
